# Weekly price-data update: a new "Femacal de La Calera - Albahaca" record
# (week of 2022-01-27, serial 44588) is inserted as row 101, pushing the
# existing rows 101:128 down to 102:129 (dimension grows from R128 to R129).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 101, shifting 101:128 -> 102:129.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly record.
$ws.Cells.Item(101, 1).Value = 3
$ws.Cells.Item(101, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(101, 3).Value = "Coquimbo"
$ws.Cells.Item(101, 4).Value = 44588
$ws.Cells.Item(101, 5).Value = 5
$ws.Cells.Item(101, 6).Value = 100112052
$ws.Cells.Item(101, 7).Value = "Albahaca"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 170
$ws.Cells.Item(101, 11).Value = 4000
$ws.Cells.Item(101, 12).Value = 4500
$ws.Cells.Item(101, 13).Value = 4235
$ws.Cells.Item(101, 14).Value = "`$/docena de matas"
$ws.Cells.Item(101, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(101, 16).Value = 706
$ws.Cells.Item(101, 17).Value = 6
$ws.Cells.Item(101, 18).Value = "Hortaliza"
